# "Generate Report for Handoff"
#
# The localization-status workbook tracks hand-off state for source files
# across locales. The file "d223aaf6-c3d8-44da-8862-95c531ca8d0f.md" (row 3
# on every sheet) has just been packaged for hand-off, so its status moves
# from "In Translation" to "Ready for handoff", its translation priority
# flips from "ht" (human translation) to "mt" (machine translation), and the
# hand-off timestamps are refreshed. The "71c82d2e..." row (row 2) is
# untouched.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the d223aaf6 file ---------------------------
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-28 02:14:47"

# --- zh-cn sheet: row 3 is the d223aaf6 file -------------------------------
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-28 02:14:43"

# --- de-de sheet: row 3 is the d223aaf6 file -------------------------------
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-28 02:14:47"

# --- Widen the "Status" columns slightly so "Ready for handoff" fits -------
$overview.Columns.Item(5).ColumnWidth = 17
$overview.Columns.Item(6).ColumnWidth = 17
$zhcn.Columns.Item(3).ColumnWidth = 17
$dede.Columns.Item(3).ColumnWidth = 17
